$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to be treated as Text so that numeric-looking
# strings (e.g. "1.0000", "13.97") are preserved verbatim instead of
# being auto-converted into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.213.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.816.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.92%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4434'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3699'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.90%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.59'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07695'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.18%  '

$ws.Range("E11").Value = '  +0.61%  '

$ws.Range("E12").Value = '  -0.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.265'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.570'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.852.38'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.59'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.58%  '

$ws.Range("E18").Value = '  +1.94%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06594'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.67%  '

$ws.Range("E20").Value = '  -0.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.53%  '

$ws.Range("E22").Value = '  +2.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.287.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.34%  '

$ws.Range("E24").Value = '  +3.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.108'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -12.51%  '

$ws.Range("E26").Value = '  +3.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.040.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.62%  '

$ws.Range("E29").Value = '  -1.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.51'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.207'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.862'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09205'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.661'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.99%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.08'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02351'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.16%  '

$ws.Range("E37").Value = '  +1.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06211'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6567'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.151'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.198'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.155'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.46%  '

$ws.Range("E43").Value = '  -0.26%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.75%  '

$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.393'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6072'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.767'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.69%  '

$ws.Range("E49").Value = '  +5.14%  '

$ws.Range("E50").Value = '  +5.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06980'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.17%  '
